$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build a 2D array with the new statistics data for columns O:AK, rows 1:13
$data = New-Object "object[,]" 13,23

$data[0,0] = 50006655
$data[0,1] = 50006724
$data[0,2] = 50006734
$data[0,3] = 50006728
$data[0,4] = 50006646
$data[0,5] = 50006662
$data[0,6] = 50006673
$data[0,7] = 50006747
$data[0,8] = 50006679
$data[0,9] = 50006678
$data[0,10] = 50008388
$data[0,11] = 50006683
$data[0,12] = 50006684
$data[0,13] = 50006688
$data[0,14] = 50006686
$data[0,15] = 50006689
$data[0,16] = 50006690
$data[0,17] = 50006691
$data[0,18] = 50006692
$data[0,19] = 50006693
$data[0,20] = 50006717
$data[0,21] = 50006694
$data[0,22] = 50006718

$data[1,0] = 0.0006590631113256104
$data[1,1] = 0.0008118840471881424
$data[1,2] = 0.0008070550108666846
$data[1,3] = 0.0008570763251744937
$data[1,4] = 0.0008655235004659501
$data[1,5] = 0.000561441236807912
$data[1,6] = 0.001168066248186017
$data[1,7] = 0.0007849240725123793
$data[1,8] = 0.0003936358679061687
$data[1,9] = 0.0003583556831178025
$data[1,10] = 0.1286282558608123
$data[1,11] = 0.06891221997728465
$data[1,12] = 0.0006809771141027355
$data[1,13] = 0.0004148715519790517
$data[1,14] = 0.08085566333922777
$data[1,15] = 0.0004540267424639727
$data[1,16] = 0.0384195839811241
$data[1,17] = 0.003252744785144668
$data[1,18] = 0.002563953913984552
$data[1,19] = 0.05267400917962285
$data[1,20] = 0.0004678051879116139
$data[1,21] = 0.0005669369624088798
$data[1,22] = 0.001488789245870938

$data[2,0] = 0.0004310556804270473
$data[2,1] = 0.0005490465730847618
$data[2,2] = 0.0006724985220634334
$data[2,3] = 0.0002837528589289472
$data[2,4] = 0.000795450890972927
$data[2,5] = 0.0004544500190832412
$data[2,6] = 0.000882207848206574
$data[2,7] = 0.0005061269696891146
$data[2,8] = 0.0003068846239756017
$data[2,9] = 0.0002598535606100043
$data[2,10] = 0.05580976789156387
$data[2,11] = 0.04409417996588735
$data[2,12] = 0.0005625292487322226
$data[2,13] = 0.0002899960030893059
$data[2,14] = 0.04806419628028925
$data[2,15] = 0.0003707960583237274
$data[2,16] = 0.0231626024143352
$data[2,17] = 0.003458655494121651
$data[2,18] = 0.00214031642389393
$data[2,19] = 0.06208008579484612
$data[2,20] = 0.0003523676149607549
$data[2,21] = 0.0005274267286427312
$data[2,22] = 0.0009629098779510806

$data[3,0] = 0.08906881885410835
$data[3,1] = 0.07090900410281664
$data[3,2] = 0.07991448230803731
$data[3,3] = 0.1253947949008939
$data[3,4] = 0.05476746493601733
$data[3,5] = 0.08032566078714996
$data[3,6] = 0.09256466860847369
$data[3,7] = 0.09094013672618315
$data[3,8] = 0.0503535194181423
$data[3,9] = 0.05357968209962193
$data[3,10] = 0.2492172796459826
$data[3,11] = 0.1742308877032095
$data[3,12] = 0.04482669135821815
$data[3,13] = 0.06541322814969423
$data[3,14] = 0.159526037009458
$data[3,15] = 0.05902582665658421
$data[3,16] = 0.1483587030267296
$data[3,17] = 0.04658401676154622
$data[3,18] = 0.07398886095962995
$data[3,19] = 0.15049717585701
$data[3,20] = 0.09582250855723605
$data[3,21] = 0.06135370622075578
$data[3,22] = 0.08966985611975332

$data[4,0] = 0.05825484639667947
$data[4,1] = 0.04795308620528023
$data[4,2] = 0.0665907162708851
$data[4,3] = 0.04151454252419538
$data[4,4] = 0.05033350192829036
$data[4,5] = 0.0650183771415485
$data[4,6] = 0.069911511645726
$data[4,7] = 0.05863911865642383
$data[4,8] = 0.03925638421793472
$data[4,9] = 0.0388521009316936
$data[4,10] = 0.1081314400053734
$data[4,11] = 0.1114833932288644
$data[4,12] = 0.03702962183408327
$data[4,13] = 0.04572397076176975
$data[4,14] = 0.09482935935446507
$data[4,15] = 0.04820540689912716
$data[4,16] = 0.08944328118188079
$data[4,17] = 0.04953295636546877
$data[4,18] = 0.06176381464322002
$data[4,19] = 0.1773716816812865
$data[4,20] = 0.07217694389111656
$data[4,21] = 0.05707792348663694
$data[4,22] = 0.05799611358802611

$data[5,0] = 0.2375344594988865
$data[5,1] = 0.2102346178435296
$data[5,2] = 0.2040233108847194
$data[5,3] = 0.2833872256230967
$data[5,4] = 0.175292348603765
$data[5,5] = 0.1970314921242014
$data[5,6] = 0.2259709133477814
$data[5,7] = 0.2304059860018139
$data[5,8] = 0.1971037343679178
$data[5,9] = 0.2055027429279748
$data[5,10] = 0.2713252696009744
$data[5,11] = 0.2041362095851902
$data[5,12] = 0.1542012756907734
$data[5,13] = 0.1949950892926545
$data[5,14] = 0.1813789336994331
$data[5,15] = 0.1958446889807068
$data[5,16] = 0.2349505508807628
$data[5,17] = 0.1492779596207345
$data[5,18] = 0.1853384738927238
$data[5,19] = 0.1580220433410729
$data[5,20] = 0.2271009213610437
$data[5,21] = 0.1881448785082401
$data[5,22] = 0.232772439511396

$data[6,0] = 0.1553577742474761
$data[6,1] = 0.1421737462024864
$data[6,2] = 0.1700074631704713
$data[6,3] = 0.09382120715809657
$data[6,4] = 0.1611007151192739
$data[6,5] = 0.1594841267181054
$data[6,6] = 0.1706695262630968
$data[6,7] = 0.1485681068744282
$data[6,8] = 0.1536651264211114
$data[6,9] = 0.1490156902225051
$data[6,10] = 0.1177237475406031
$data[6,11] = 0.1306186154788007
$data[6,12] = 0.1273797987795441
$data[6,13] = 0.1363019379062351
$data[6,14] = 0.1078196914155993
$data[6,15] = 0.1599430868842373
$data[6,16] = 0.141648233352914
$data[6,17] = 0.1587278035311832
$data[6,18] = 0.1547153314606821
$data[6,19] = 0.186240276008566
$data[6,20] = 0.1710605441821229
$data[6,21] = 0.1750329302235859
$data[6,22] = 0.1505511152380557

$data[7,0] = 0.1580091813912946
$data[7,1] = 0.1678594339425228
$data[7,2] = 0.1528343926443516
$data[7,3] = 0.1761427307327432
$data[7,4] = 0.1556812892236656
$data[7,5] = 0.1554459859901691
$data[7,6] = 0.1397998883586611
$data[7,7] = 0.151870584993588
$data[7,8] = 0.1658222088116965
$data[7,9] = 0.1726480946537363
$data[7,10] = 0.04019632317337211
$data[7,11] = 0.1495265073417985
$data[7,12] = 0.1787529173936202
$data[7,13] = 0.1851761700698496
$data[7,14] = 0.1638966163474531
$data[7,15] = 0.1513482860325678
$data[7,16] = 0.1554515547221037
$data[7,17] = 0.1496264695937122
$data[7,18] = 0.178379106867289
$data[7,19] = 0.06772372834554428
$data[7,20] = 0.1491915203082791
$data[7,21] = 0.1584778742187673
$data[7,22] = 0.1672358479890854

$data[8,0] = 0.1033448148256246
$data[8,1] = 0.1135170068746685
$data[8,2] = 0.1273530326804043
$data[8,3] = 0.05831569716360082
$data[8,4] = 0.1430773631843755
$data[8,5] = 0.1258233750361567
$data[8,6] = 0.1055869552604115
$data[8,7] = 0.09792768709677399
$data[8,8] = 0.1292775642338352
$data[8,9] = 0.1251918812560272
$data[8,10] = 0.01744054952302002
$data[8,11] = 0.09567604594037389
$data[8,12] = 0.147660974572717
$data[8,13] = 0.129438494713548
$data[8,14] = 0.09742742576669193
$data[8,15] = 0.1236036176864214
$data[8,16] = 0.09371945720410192
$data[8,17] = 0.1590983754673216
$data[8,18] = 0.1489058481219992
$data[8,19] = 0.07981725582537753
$data[8,20] = 0.1123763941526229
$data[8,21] = 0.1474334402299494
$data[8,22] = 0.1081637649001231

$data[9,0] = 0.0851011145029768
$data[9,1] = 0.1104032799651282
$data[9,2] = 0.08381528559816662
$data[9,3] = 0.1266875224080188
$data[9,4] = 0.1005949324470683
$data[9,5] = 0.09141195294190554
$data[9,6] = 0.0817241207882579
$data[9,7] = 0.09513110640856311
$data[9,8] = 0.1064148870383514
$data[9,9] = 0.1069774120706555
$data[9,10] = 0.008039259208169031
$data[9,11] = 0.01300229256056102
$data[9,12] = 0.1212962956546219
$data[9,13] = 0.1035824834047581
$data[9,14] = 0.04152046743516594
$data[9,15] = 0.1118463730962214
$data[9,16] = 0.04669456587819966
$data[9,17] = 0.1036232878244818
$data[9,18] = 0.07838425553723283
$data[9,19] = 0.03009942252963853
$data[9,20] = 0.08132047960169961
$data[9,21] = 0.08499693237849376
$data[9,22] = 0.09313874739957191

$data[10,0] = 0.0556597967429819
$data[10,1] = 0.07466157603676152
$data[10,2] = 0.06984115696222662
$data[10,3] = 0.04194252672488846
$data[10,4] = 0.09245078683514017
$data[10,5] = 0.07399200670594544
$data[10,6] = 0.0617239483284012
$data[10,7] = 0.06134149823641623
$data[10,8] = 0.08296269536584784
$data[10,9] = 0.07757226337126422
$data[10,10] = 0.003488107550129054
$data[10,11] = 0.008319648218029829
$data[10,12] = 0.1001982484514187
$data[10,13] = 0.07240435270664437
$data[10,14] = 0.02468160935221676
$data[10,15] = 0.09134306507325334
$data[10,16] = 0.02815146735784817
$data[10,17] = 0.1101830230855795
$data[10,18] = 0.06543296608647582
$data[10,19] = 0.03547432143703521
$data[10,20] = 0.0612534965091699
$data[10,21] = 0.07907343666318385
$data[10,22] = 0.06023970158286124

$data[11,0] = 0.03420652989527234
$data[11,1] = 0.03634712670357927
$data[11,2] = 0.02407746419443166
$data[11,3] = 0.03880554472453106
$data[11,4] = 0.03389226725472702
$data[11,5] = 0.02788226140403236
$data[11,6] = 0.0284845827053239
$data[11,7] = 0.03884018219708657
$data[11,8] = 0.04183114815154228
$data[11,9] = 0.04060101251808564
$data[11,10] = 0
$data[11,11] = 0
$data[11,12] = 0.04786840462646322
$data[11,13] = 0.0389990178585309
$data[11,14] = 0
$data[11,15] = 0.03193447415413185
$data[11,16] = 0
$data[11,17] = 0.03229515275740245
$data[11,18] = 0.02637225626253302
$data[11,19] = 0
$data[11,20] = 0.01647069871253892
$data[11,21] = 0.02451136425615071
$data[11,22] = 0.02294226113343875

$data[12,0] = 0.0223725448529466
$data[12,1] = 0.02458019150295394
$data[12,2] = 0.02006314175337578
$data[12,3] = 0.01284737885583181
$data[12,4] = 0.03114835607623777
$data[12,5] = 0.02256886989489451
$data[12,6] = 0.02151361059747367
$data[12,7] = 0.02504454176652161
$data[12,8] = 0.03261221148173871
$data[12,9] = 0.02944091070470815
$data[12,10] = 0
$data[12,11] = 0
$data[12,12] = 0.03954226527570499
$data[12,13] = 0.02726038758124703
$data[12,14] = 0
$data[12,15] = 0.02608035173596111
$data[12,16] = 0
$data[12,17] = 0.03433955471330351
$data[12,18] = 0.02201481583033596
$data[12,19] = 0
$data[12,20] = 0.01240631992129822
$data[12,21] = 0.02280315012318459
$data[12,22] = 0.01483845341386649

# Write the whole block in one shot
$ws.Range("O1:AK13").Value = $data

# Copy header row (N1) formatting (bold font + border + centered/top alignment) to the new header cells O1:AK1
$ws.Range("N1").Copy() | Out-Null
$ws.Range("O1:AK1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
